$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values scraped on this run (Price/Volume columns); some Price values look
# like plain numbers (e.g. "1.00", "0.0370") so we force Text format before writing
# them to keep the exact display (trailing/leading zeros), then restore the default
# "Normal" style so no stray formatting is left behind.

$ws.Range("D2").Value = "63.686.05"
$ws.Range("E2").Value = "  -1.07%  "
$ws.Range("D3").Value = "3.131.62"
$ws.Range("E3").Value = "  -0.96%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.73%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "3.130.74"
$ws.Range("E8").Value = "  -0.81%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.528"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.43%  "
$ws.Range("E10").Value = "  +1.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.74"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.457"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.57%  "
$ws.Range("E13").Value = "  -3.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.75"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.76%  "
$ws.Range("E15").Value = "  -1.82%  "
$ws.Range("D16").Value = "3.653.84"
$ws.Range("E16").Value = "  -0.83%  "
$ws.Range("D17").Value = "63.577.18"
$ws.Range("E17").Value = "  -0.96%  "
$ws.Range("D18").Value = "3.135.96"
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.05"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "463.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.62%  "
$ws.Range("E22").Value = "  -0.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.41"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.70%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.19"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.53%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.22"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.81%  "
$ws.Range("E29").Value = "  -1.36%  "
$ws.Range("B30").Value = "FirstDigitalUSD"
$ws.Range("C30").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.15%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.21"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.96"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.89"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.91%  "
$ws.Range("E34").Value = "  -0.21%  "
$ws.Range("D35").Value = "0.0₃0840"
$ws.Range("E35").Value = "  -5.63%  "
$ws.Range("E36").Value = "  -1.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.28"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.68%  "
$ws.Range("E38").Value = "  -3.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "51.08"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "438.06"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.58%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.77"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.23%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.907.95"
$ws.Range("E43").Value = "  -1.38%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0370"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.65%  "
$ws.Range("E45").Value = "  -2.55%  "
$ws.Range("E46").Value = "  -4.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "36.97"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "126.50"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.35%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("E50").Value = "  -1.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.12"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.01%  "
